$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 13: was "519033 - Carlos Yujiro Shigue" in B/C only (no A cell).
# Becomes "Programa resumido:" / "Semestral" / "Semestral", height 60.
# ------------------------------------------------------------------
$ws.Cells.Item(13,1).Value = "Programa resumido:"
$ws.Cells.Item(13,2).Value = "Semestral"
$ws.Cells.Item(13,3).Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# ------------------------------------------------------------------
# Row 14: keep A only ("Short syllabus:"), drop B/C (used to hold the
# long "Metodologia..." text). Height stays 60.
# ------------------------------------------------------------------
$ws.Cells.Item(14,1).Value = "Short syllabus:"
$ws.Cells.Item(14,2).Clear()
$ws.Cells.Item(14,3).Clear()
$ws.Rows.Item(14).RowHeight = 60

# ------------------------------------------------------------------
# Row 15: "Programa:" plus new B/C = "01/01/2012". Height 60 -> 120.
# B15 is a brand-new cell: copy number/format (not value) from a
# known-good column-B cell first so it lands on style id 2 instead of
# picking up the overlapping (wrong) column style.
# ------------------------------------------------------------------
$ws.Cells.Item(19,2).Copy()
$ws.Cells.Item(15,2).PasteSpecial(-4122)
$ws.Cells.Item(19,3).Copy()
$ws.Cells.Item(15,3).PasteSpecial(-4122)

$ws.Cells.Item(15,1).Value = "Programa:"
$ws.Cells.Item(15,2).Value = "01/01/2012"
$ws.Cells.Item(15,3).Value = "01/01/2012"
$ws.Rows.Item(15).RowHeight = 120

# ------------------------------------------------------------------
# Row 16: keep A only ("Syllabus:"), drop B/C (used to hold
# "Apresentação de projetos..."). Height stays 120.
# ------------------------------------------------------------------
$ws.Cells.Item(16,1).Value = "Syllabus:"
$ws.Cells.Item(16,2).Clear()
$ws.Cells.Item(16,3).Clear()
$ws.Rows.Item(16).RowHeight = 120

# ------------------------------------------------------------------
# Row 17: "Avaliação:" only, height goes back to default (was 120).
# ------------------------------------------------------------------
$ws.Cells.Item(17,1).Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# ------------------------------------------------------------------
# Row 18: "Método:" plus new B/C = "519033 - Carlos Yujiro Shigue".
# Height none -> 60. Fix B18's style the same way as B15 above.
# ------------------------------------------------------------------
$ws.Cells.Item(19,2).Copy()
$ws.Cells.Item(18,2).PasteSpecial(-4122)
$ws.Cells.Item(19,3).Copy()
$ws.Cells.Item(18,3).PasteSpecial(-4122)

$ws.Cells.Item(18,1).Value = "Método:"
$ws.Cells.Item(18,2).Value = "519033 - Carlos Yujiro Shigue"
$ws.Cells.Item(18,3).Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(18).RowHeight = 60

# ------------------------------------------------------------------
# Row 19: "Critério:" plus the "Aulas expositivas..." text. Height
# stays 60 (cells already existed here, just change A's text).
# ------------------------------------------------------------------
$ws.Cells.Item(19,1).Value = "Critério:"
$ws.Cells.Item(19,2).Value = "Aulas expositivas, reuniões com supervisor, desenvolvimento e elaboração de projeto."
$ws.Cells.Item(19,3).Value = "Aulas expositivas, reuniões com supervisor, desenvolvimento e elaboração de projeto."
$ws.Rows.Item(19).RowHeight = 60

# ------------------------------------------------------------------
# Row 20: "Norma de recuperação:" plus "Nota de avaliação do
# projeto.". Height stays 60.
# ------------------------------------------------------------------
$ws.Cells.Item(20,1).Value = "Norma de recuperação:"
$ws.Cells.Item(20,2).Value = "Nota de avaliação do projeto."
$ws.Cells.Item(20,3).Value = "Nota de avaliação do projeto."
$ws.Rows.Item(20).RowHeight = 60

# ------------------------------------------------------------------
# Row 21: "Bibliografia:" plus "A critério da Comissão de Curso...".
# Height 60 -> 120.
# ------------------------------------------------------------------
$ws.Cells.Item(21,1).Value = "Bibliografia:"
$ws.Cells.Item(21,2).Value = "A critério da Comissão de Curso poderá ser oferecida recuperação."
$ws.Cells.Item(21,3).Value = "A critério da Comissão de Curso poderá ser oferecida recuperação."
$ws.Rows.Item(21).RowHeight = 120

# ------------------------------------------------------------------
# Row 22: keep A only ("Requisitos:"), drop B/C (used to hold the
# long bibliography text). Height goes back to default (was 120).
# ------------------------------------------------------------------
$ws.Cells.Item(22,1).Value = "Requisitos:"
$ws.Cells.Item(22,2).Clear()
$ws.Cells.Item(22,3).Clear()
$ws.Rows.Item(22).AutoFit()

# ------------------------------------------------------------------
# Row 23: drop A (used to hold "Requisitos:"), add new B/C with the
# first prerequisite line. Height none -> 30. Fix B23's style.
# ------------------------------------------------------------------
$ws.Cells.Item(23,1).Clear()

$ws.Cells.Item(19,2).Copy()
$ws.Cells.Item(23,2).PasteSpecial(-4122)
$ws.Cells.Item(19,3).Copy()
$ws.Cells.Item(23,3).PasteSpecial(-4122)

$ws.Cells.Item(23,2).Value = "LOQ4050 -  Engenharia Econômica  (Requisito)`n"
$ws.Cells.Item(23,3).Value = "LOQ4050 -  Engenharia Econômica  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# ------------------------------------------------------------------
# Row 24: second prerequisite line moves here from row 25.
# Height stays 30.
# ------------------------------------------------------------------
$ws.Cells.Item(24,2).Value = "LOQ4234 -  Empreendedorismo  (Requisito)`n"
$ws.Cells.Item(24,3).Value = "LOQ4234 -  Empreendedorismo  (Requisito)`n"
$ws.Rows.Item(24).RowHeight = 30

# ------------------------------------------------------------------
# Row 25 no longer exists in the final layout: clear its cells and
# let AutoFit drop the now-empty row entirely.
# ------------------------------------------------------------------
$ws.Cells.Item(25,2).Clear()
$ws.Cells.Item(25,3).Clear()
$ws.Rows.Item(25).AutoFit()
